$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Cells.Item(20, 2).Value2 = 6839250
$ws.Cells.Item(20, 5).Value2 = 'Dungannon Swifts'
$ws.Cells.Item(20, 6).Value2 = 'Loughgall'
$ws.Cells.Item(20, 7).Value2 = 1
$ws.Cells.Item(20, 8).Value2 = 2
$ws.Cells.Item(20, 9).Value2 = 0
$ws.Cells.Item(20, 10).Value2 = 1
$ws.Cells.Item(20, 11).Value2 = 'A'
$ws.Cells.Item(20, 12).Value2 = 2.3
$ws.Cells.Item(20, 13).Value2 = 3.6
$ws.Cells.Item(20, 14).Value2 = 2.625
$ws.Cells.Item(20, 15).Value2 = 2.25
$ws.Cells.Item(20, 16).Value2 = 3.5
$ws.Cells.Item(20, 17).Value2 = 2.7
$ws.Cells.Item(20, 18).Value2 = -0.25
$ws.Cells.Item(20, 19).Value2 = 2.05
$ws.Cells.Item(20, 20).Value2 = 1.75
$ws.Cells.Item(20, 21).Value2 = 2.5
$ws.Cells.Item(20, 22).Value2 = 1.95
$ws.Cells.Item(20, 23).Value2 = 1.85
$ws.Cells.Item(20, 24).Value2 = -1
$ws.Cells.Item(20, 25).Value2 = -1
$ws.Cells.Item(20, 26).Value2 = 1.7
$ws.Cells.Item(20, 27).Value2 = -1
$ws.Cells.Item(20, 28).Value2 = 0.75
$ws.Cells.Item(20, 29).Value2 = 0.95
$ws.Cells.Item(20, 30).Value2 = -1

# Row 21
$ws.Cells.Item(21, 2).Value2 = 6840306
$ws.Cells.Item(21, 5).Value2 = 'Glenavon'
$ws.Cells.Item(21, 6).Value2 = 'Newry City'
$ws.Cells.Item(21, 7).Value2 = 1
$ws.Cells.Item(21, 8).Value2 = 3
$ws.Cells.Item(21, 9).Value2 = 0
$ws.Cells.Item(21, 10).Value2 = 1
$ws.Cells.Item(21, 11).Value2 = 'A'
$ws.Cells.Item(21, 12).Value2 = 1.727
$ws.Cells.Item(21, 13).Value2 = 4
$ws.Cells.Item(21, 14).Value2 = 3.8
$ws.Cells.Item(21, 15).Value2 = 1.65
$ws.Cells.Item(21, 16).Value2 = 4
$ws.Cells.Item(21, 17).Value2 = 4.2
$ws.Cells.Item(21, 18).Value2 = -0.75
$ws.Cells.Item(21, 19).Value2 = 1.825
$ws.Cells.Item(21, 20).Value2 = 1.975
$ws.Cells.Item(21, 21).Value2 = 2.75
$ws.Cells.Item(21, 22).Value2 = 1.975
$ws.Cells.Item(21, 23).Value2 = 1.825
$ws.Cells.Item(21, 24).Value2 = -1
$ws.Cells.Item(21, 25).Value2 = -1
$ws.Cells.Item(21, 26).Value2 = 3.2
$ws.Cells.Item(21, 27).Value2 = -1
$ws.Cells.Item(21, 28).Value2 = 0.9750000000000001
$ws.Cells.Item(21, 29).Value2 = 0.9750000000000001
$ws.Cells.Item(21, 30).Value2 = -1

# Row 22
$ws.Cells.Item(22, 2).Value2 = 6840432
$ws.Cells.Item(22, 5).Value2 = 'Cliftonville'
$ws.Cells.Item(22, 6).Value2 = 'Carrick Rangers'
$ws.Cells.Item(22, 7).Value2 = 4
$ws.Cells.Item(22, 8).Value2 = 0
$ws.Cells.Item(22, 9).Value2 = 3
$ws.Cells.Item(22, 10).Value2 = 0
$ws.Cells.Item(22, 11).Value2 = 'H'
$ws.Cells.Item(22, 12).Value2 = 1.4
$ws.Cells.Item(22, 13).Value2 = 4.5
$ws.Cells.Item(22, 14).Value2 = 6.5
$ws.Cells.Item(22, 15).Value2 = 1.533
$ws.Cells.Item(22, 16).Value2 = 4.2
$ws.Cells.Item(22, 17).Value2 = 5
$ws.Cells.Item(22, 18).Value2 = -1
$ws.Cells.Item(22, 19).Value2 = 1.85
$ws.Cells.Item(22, 20).Value2 = 1.95
$ws.Cells.Item(22, 21).Value2 = 2.75
$ws.Cells.Item(22, 22).Value2 = 1.85
$ws.Cells.Item(22, 23).Value2 = 1.95
$ws.Cells.Item(22, 24).Value2 = 0.5329999999999999
$ws.Cells.Item(22, 25).Value2 = -1
$ws.Cells.Item(22, 26).Value2 = -1
$ws.Cells.Item(22, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(22, 28).Value2 = -1
$ws.Cells.Item(22, 29).Value2 = 0.8500000000000001
$ws.Cells.Item(22, 30).Value2 = -1

# Row 23
$ws.Cells.Item(23, 2).Value2 = 6973551
$ws.Cells.Item(23, 5).Value2 = 'Linfield'
$ws.Cells.Item(23, 6).Value2 = 'Cliftonville'
$ws.Cells.Item(23, 7).Value2 = 2
$ws.Cells.Item(23, 8).Value2 = 1
$ws.Cells.Item(23, 9).Value2 = 1
$ws.Cells.Item(23, 10).Value2 = 0
$ws.Cells.Item(23, 11).Value2 = 'H'
$ws.Cells.Item(23, 12).Value2 = 1.727
$ws.Cells.Item(23, 13).Value2 = 3.75
$ws.Cells.Item(23, 14).Value2 = 4
$ws.Cells.Item(23, 15).Value2 = 1.8
$ws.Cells.Item(23, 16).Value2 = 3.6
$ws.Cells.Item(23, 17).Value2 = 3.8
$ws.Cells.Item(23, 18).Value2 = -0.5
$ws.Cells.Item(23, 19).Value2 = 1.825
$ws.Cells.Item(23, 20).Value2 = 1.975
$ws.Cells.Item(23, 21).Value2 = 2.5
$ws.Cells.Item(23, 22).Value2 = 1.875
$ws.Cells.Item(23, 23).Value2 = 1.925
$ws.Cells.Item(23, 24).Value2 = 0.8
$ws.Cells.Item(23, 25).Value2 = -1
$ws.Cells.Item(23, 26).Value2 = -1
$ws.Cells.Item(23, 27).Value2 = 0.825
$ws.Cells.Item(23, 28).Value2 = -1
$ws.Cells.Item(23, 29).Value2 = 0.875
$ws.Cells.Item(23, 30).Value2 = -1

# Row 24
$ws.Cells.Item(24, 2).Value2 = 6978093
$ws.Cells.Item(24, 5).Value2 = 'Crusaders'
$ws.Cells.Item(24, 6).Value2 = 'Carrick Rangers'
$ws.Cells.Item(24, 7).Value2 = 9
$ws.Cells.Item(24, 8).Value2 = 0
$ws.Cells.Item(24, 9).Value2 = 5
$ws.Cells.Item(24, 10).Value2 = 0
$ws.Cells.Item(24, 11).Value2 = 'H'
$ws.Cells.Item(24, 12).Value2 = 1.333
$ws.Cells.Item(24, 13).Value2 = 5
$ws.Cells.Item(24, 14).Value2 = 7
$ws.Cells.Item(24, 15).Value2 = 1.4
$ws.Cells.Item(24, 16).Value2 = 4.2
$ws.Cells.Item(24, 17).Value2 = 6.5
$ws.Cells.Item(24, 18).Value2 = -1.25
$ws.Cells.Item(24, 19).Value2 = 1.85
$ws.Cells.Item(24, 20).Value2 = 1.95
$ws.Cells.Item(24, 21).Value2 = 2.75
$ws.Cells.Item(24, 22).Value2 = 1.775
$ws.Cells.Item(24, 23).Value2 = 2.025
$ws.Cells.Item(24, 24).Value2 = 0.3999999999999999
$ws.Cells.Item(24, 25).Value2 = -1
$ws.Cells.Item(24, 26).Value2 = -1
$ws.Cells.Item(24, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(24, 28).Value2 = -1
$ws.Cells.Item(24, 29).Value2 = 0.7749999999999999
$ws.Cells.Item(24, 30).Value2 = -1

# Row 47
$ws.Cells.Item(47, 2).Value2 = 6840444
$ws.Cells.Item(47, 5).Value2 = 'Larne FC'
$ws.Cells.Item(47, 6).Value2 = 'Carrick Rangers'
$ws.Cells.Item(47, 7).Value2 = 4
$ws.Cells.Item(47, 8).Value2 = 1
$ws.Cells.Item(47, 9).Value2 = 3
$ws.Cells.Item(47, 10).Value2 = 0
$ws.Cells.Item(47, 11).Value2 = 'H'
$ws.Cells.Item(47, 12).Value2 = 1.25
$ws.Cells.Item(47, 13).Value2 = 6
$ws.Cells.Item(47, 14).Value2 = 8
$ws.Cells.Item(47, 15).Value2 = 1.25
$ws.Cells.Item(47, 16).Value2 = 6
$ws.Cells.Item(47, 17).Value2 = 7.5
$ws.Cells.Item(47, 18).Value2 = -1.75
$ws.Cells.Item(47, 19).Value2 = 1.975
$ws.Cells.Item(47, 20).Value2 = 1.825
$ws.Cells.Item(47, 21).Value2 = 3
$ws.Cells.Item(47, 22).Value2 = 1.9
$ws.Cells.Item(47, 23).Value2 = 1.9
$ws.Cells.Item(47, 24).Value2 = 0.25
$ws.Cells.Item(47, 25).Value2 = -1
$ws.Cells.Item(47, 26).Value2 = -1
$ws.Cells.Item(47, 27).Value2 = 0.9750000000000001
$ws.Cells.Item(47, 28).Value2 = -1
$ws.Cells.Item(47, 29).Value2 = 0.8999999999999999
$ws.Cells.Item(47, 30).Value2 = -1

# Row 48
$ws.Cells.Item(48, 2).Value2 = 6837581
$ws.Cells.Item(48, 5).Value2 = 'Glentoran'
$ws.Cells.Item(48, 6).Value2 = 'Coleraine'
$ws.Cells.Item(48, 7).Value2 = 1
$ws.Cells.Item(48, 8).Value2 = 2
$ws.Cells.Item(48, 9).Value2 = 1
$ws.Cells.Item(48, 10).Value2 = 0
$ws.Cells.Item(48, 11).Value2 = 'A'
$ws.Cells.Item(48, 12).Value2 = 1.666
$ws.Cells.Item(48, 13).Value2 = 3.75
$ws.Cells.Item(48, 14).Value2 = 4.5
$ws.Cells.Item(48, 15).Value2 = 1.7
$ws.Cells.Item(48, 16).Value2 = 3.75
$ws.Cells.Item(48, 17).Value2 = 4.2
$ws.Cells.Item(48, 18).Value2 = -0.75
$ws.Cells.Item(48, 19).Value2 = 1.95
$ws.Cells.Item(48, 20).Value2 = 1.85
$ws.Cells.Item(48, 21).Value2 = 2.5
$ws.Cells.Item(48, 22).Value2 = 1.8
$ws.Cells.Item(48, 23).Value2 = 2
$ws.Cells.Item(48, 24).Value2 = -1
$ws.Cells.Item(48, 25).Value2 = -1
$ws.Cells.Item(48, 26).Value2 = 3.2
$ws.Cells.Item(48, 27).Value2 = -1
$ws.Cells.Item(48, 28).Value2 = 0.8500000000000001
$ws.Cells.Item(48, 29).Value2 = 0.8
$ws.Cells.Item(48, 30).Value2 = -1

# Row 49
$ws.Cells.Item(49, 2).Value2 = 6840445
$ws.Cells.Item(49, 5).Value2 = 'Newry City'
$ws.Cells.Item(49, 6).Value2 = 'Dungannon Swifts'
$ws.Cells.Item(49, 7).Value2 = 2
$ws.Cells.Item(49, 8).Value2 = 2
$ws.Cells.Item(49, 9).Value2 = 1
$ws.Cells.Item(49, 10).Value2 = 2
$ws.Cells.Item(49, 11).Value2 = 'D'
$ws.Cells.Item(49, 12).Value2 = 2.625
$ws.Cells.Item(49, 13).Value2 = 3.4
$ws.Cells.Item(49, 14).Value2 = 2.375
$ws.Cells.Item(49, 15).Value2 = 2.8
$ws.Cells.Item(49, 16).Value2 = 3.4
$ws.Cells.Item(49, 17).Value2 = 2.3
$ws.Cells.Item(49, 18).Value2 = 0.25
$ws.Cells.Item(49, 19).Value2 = 1.775
$ws.Cells.Item(49, 20).Value2 = 2.025
$ws.Cells.Item(49, 21).Value2 = 2.75
$ws.Cells.Item(49, 22).Value2 = 1.975
$ws.Cells.Item(49, 23).Value2 = 1.825
$ws.Cells.Item(49, 24).Value2 = -1
$ws.Cells.Item(49, 25).Value2 = 2.4
$ws.Cells.Item(49, 26).Value2 = -1
$ws.Cells.Item(49, 27).Value2 = 0.3875
$ws.Cells.Item(49, 28).Value2 = -0.5
$ws.Cells.Item(49, 29).Value2 = 0.9750000000000001
$ws.Cells.Item(49, 30).Value2 = -1

# Row 64
$ws.Cells.Item(64, 2).Value2 = 6840450
$ws.Cells.Item(64, 5).Value2 = 'Coleraine'
$ws.Cells.Item(64, 6).Value2 = 'Glenavon'
$ws.Cells.Item(64, 7).Value2 = 4
$ws.Cells.Item(64, 8).Value2 = 1
$ws.Cells.Item(64, 9).Value2 = 1
$ws.Cells.Item(64, 10).Value2 = 1
$ws.Cells.Item(64, 11).Value2 = 'H'
$ws.Cells.Item(64, 12).Value2 = 1.5
$ws.Cells.Item(64, 13).Value2 = 4.25
$ws.Cells.Item(64, 14).Value2 = 4.75
$ws.Cells.Item(64, 15).Value2 = 1.5
$ws.Cells.Item(64, 16).Value2 = 4.333
$ws.Cells.Item(64, 17).Value2 = 5
$ws.Cells.Item(64, 18).Value2 = -1
$ws.Cells.Item(64, 19).Value2 = 1.8
$ws.Cells.Item(64, 20).Value2 = 2
$ws.Cells.Item(64, 21).Value2 = 2.75
$ws.Cells.Item(64, 22).Value2 = 1.825
$ws.Cells.Item(64, 23).Value2 = 1.975
$ws.Cells.Item(64, 24).Value2 = 0.5
$ws.Cells.Item(64, 25).Value2 = -1
$ws.Cells.Item(64, 26).Value2 = -1
$ws.Cells.Item(64, 27).Value2 = 0.8
$ws.Cells.Item(64, 28).Value2 = -1
$ws.Cells.Item(64, 29).Value2 = 0.825
$ws.Cells.Item(64, 30).Value2 = -1

# Row 67
$ws.Cells.Item(67, 2).Value2 = 6840449
$ws.Cells.Item(67, 5).Value2 = 'Ballymena Utd'
$ws.Cells.Item(67, 6).Value2 = 'Newry City'
$ws.Cells.Item(67, 7).Value2 = 0
$ws.Cells.Item(67, 8).Value2 = 2
$ws.Cells.Item(67, 9).Value2 = 0
$ws.Cells.Item(67, 10).Value2 = 1
$ws.Cells.Item(67, 11).Value2 = 'A'
$ws.Cells.Item(67, 12).Value2 = 2
$ws.Cells.Item(67, 13).Value2 = 3.6
$ws.Cells.Item(67, 14).Value2 = 3
$ws.Cells.Item(67, 15).Value2 = 2
$ws.Cells.Item(67, 16).Value2 = 3.5
$ws.Cells.Item(67, 17).Value2 = 3.25
$ws.Cells.Item(67, 18).Value2 = -0.25
$ws.Cells.Item(67, 19).Value2 = 1.775
$ws.Cells.Item(67, 20).Value2 = 2.025
$ws.Cells.Item(67, 21).Value2 = 2.5
$ws.Cells.Item(67, 22).Value2 = 1.925
$ws.Cells.Item(67, 23).Value2 = 1.875
$ws.Cells.Item(67, 24).Value2 = -1
$ws.Cells.Item(67, 25).Value2 = -1
$ws.Cells.Item(67, 26).Value2 = 2.25
$ws.Cells.Item(67, 27).Value2 = -1
$ws.Cells.Item(67, 28).Value2 = 1.025
$ws.Cells.Item(67, 29).Value2 = -1
$ws.Cells.Item(67, 30).Value2 = 0.875

# Row 68
$ws.Cells.Item(68, 2).Value2 = 6840452
$ws.Cells.Item(68, 5).Value2 = 'Crusaders'
$ws.Cells.Item(68, 6).Value2 = 'Ballymena Utd'
$ws.Cells.Item(68, 7).Value2 = 1
$ws.Cells.Item(68, 8).Value2 = 0
$ws.Cells.Item(68, 9).Value2 = 1
$ws.Cells.Item(68, 10).Value2 = 0
$ws.Cells.Item(68, 11).Value2 = 'H'
$ws.Cells.Item(68, 12).Value2 = 1.222
$ws.Cells.Item(68, 13).Value2 = 5.25
$ws.Cells.Item(68, 14).Value2 = 10
$ws.Cells.Item(68, 15).Value2 = 1.166
$ws.Cells.Item(68, 16).Value2 = 6
$ws.Cells.Item(68, 17).Value2 = 12
$ws.Cells.Item(68, 18).Value2 = -2
$ws.Cells.Item(68, 19).Value2 = 1.85
$ws.Cells.Item(68, 20).Value2 = 1.95
$ws.Cells.Item(68, 21).Value2 = 3.25
$ws.Cells.Item(68, 22).Value2 = 1.925
$ws.Cells.Item(68, 23).Value2 = 1.875
$ws.Cells.Item(68, 24).Value2 = 0.1659999999999999
$ws.Cells.Item(68, 25).Value2 = -1
$ws.Cells.Item(68, 26).Value2 = -1
$ws.Cells.Item(68, 27).Value2 = -1
$ws.Cells.Item(68, 28).Value2 = 0.95
$ws.Cells.Item(68, 29).Value2 = -1
$ws.Cells.Item(68, 30).Value2 = 0.875

# Row 69
$ws.Cells.Item(69, 2).Value2 = 6840454
$ws.Cells.Item(69, 5).Value2 = 'Glenavon'
$ws.Cells.Item(69, 6).Value2 = 'Cliftonville'
$ws.Cells.Item(69, 7).Value2 = 0
$ws.Cells.Item(69, 8).Value2 = 1
$ws.Cells.Item(69, 9).Value2 = 0
$ws.Cells.Item(69, 10).Value2 = 0
$ws.Cells.Item(69, 11).Value2 = 'A'
$ws.Cells.Item(69, 12).Value2 = 4.75
$ws.Cells.Item(69, 13).Value2 = 4.2
$ws.Cells.Item(69, 14).Value2 = 1.533
$ws.Cells.Item(69, 15).Value2 = 7.5
$ws.Cells.Item(69, 16).Value2 = 5.5
$ws.Cells.Item(69, 17).Value2 = 1.285
$ws.Cells.Item(69, 18).Value2 = 1.5
$ws.Cells.Item(69, 19).Value2 = 1.95
$ws.Cells.Item(69, 20).Value2 = 1.85
$ws.Cells.Item(69, 21).Value2 = 3
$ws.Cells.Item(69, 22).Value2 = 1.85
$ws.Cells.Item(69, 23).Value2 = 1.95
$ws.Cells.Item(69, 24).Value2 = -1
$ws.Cells.Item(69, 25).Value2 = -1
$ws.Cells.Item(69, 26).Value2 = 0.2849999999999999
$ws.Cells.Item(69, 27).Value2 = 0.95
$ws.Cells.Item(69, 28).Value2 = -1
$ws.Cells.Item(69, 29).Value2 = -1
$ws.Cells.Item(69, 30).Value2 = 0.95

# Row 70
$ws.Cells.Item(70, 2).Value2 = 6840453
$ws.Cells.Item(70, 5).Value2 = 'Dungannon Swifts'
$ws.Cells.Item(70, 6).Value2 = 'Larne FC'
$ws.Cells.Item(70, 7).Value2 = 0
$ws.Cells.Item(70, 8).Value2 = 0
$ws.Cells.Item(70, 9).Value2 = 0
$ws.Cells.Item(70, 10).Value2 = 0
$ws.Cells.Item(70, 11).Value2 = 'D'
$ws.Cells.Item(70, 12).Value2 = 7
$ws.Cells.Item(70, 13).Value2 = 4.2
$ws.Cells.Item(70, 14).Value2 = 1.363
$ws.Cells.Item(70, 15).Value2 = 7
$ws.Cells.Item(70, 16).Value2 = 4.2
$ws.Cells.Item(70, 17).Value2 = 1.363
$ws.Cells.Item(70, 18).Value2 = 1.5
$ws.Cells.Item(70, 19).Value2 = 1.8
$ws.Cells.Item(70, 20).Value2 = 2
$ws.Cells.Item(70, 21).Value2 = 3
$ws.Cells.Item(70, 22).Value2 = 2
$ws.Cells.Item(70, 23).Value2 = 1.8
$ws.Cells.Item(70, 24).Value2 = -1
$ws.Cells.Item(70, 25).Value2 = 3.2
$ws.Cells.Item(70, 26).Value2 = -1
$ws.Cells.Item(70, 27).Value2 = 0.8
$ws.Cells.Item(70, 28).Value2 = -1
$ws.Cells.Item(70, 29).Value2 = -1
$ws.Cells.Item(70, 30).Value2 = 0.8

# Row 73
$ws.Cells.Item(73, 2).Value2 = 6840317
$ws.Cells.Item(73, 5).Value2 = 'Carrick Rangers'
$ws.Cells.Item(73, 6).Value2 = 'Linfield'
$ws.Cells.Item(73, 7).Value2 = 1
$ws.Cells.Item(73, 8).Value2 = 2
$ws.Cells.Item(73, 9).Value2 = 0
$ws.Cells.Item(73, 10).Value2 = 0
$ws.Cells.Item(73, 11).Value2 = 'A'
$ws.Cells.Item(73, 12).Value2 = 8
$ws.Cells.Item(73, 13).Value2 = 4.75
$ws.Cells.Item(73, 14).Value2 = 1.285
$ws.Cells.Item(73, 15).Value2 = 8
$ws.Cells.Item(73, 16).Value2 = 5
$ws.Cells.Item(73, 17).Value2 = 1.285
$ws.Cells.Item(73, 18).Value2 = 1.75
$ws.Cells.Item(73, 19).Value2 = 1.775
$ws.Cells.Item(73, 20).Value2 = 2.025
$ws.Cells.Item(73, 21).Value2 = 3
$ws.Cells.Item(73, 22).Value2 = 1.8
$ws.Cells.Item(73, 23).Value2 = 2
$ws.Cells.Item(73, 24).Value2 = -1
$ws.Cells.Item(73, 25).Value2 = -1
$ws.Cells.Item(73, 26).Value2 = 0.2849999999999999
$ws.Cells.Item(73, 27).Value2 = 0.7749999999999999
$ws.Cells.Item(73, 28).Value2 = -1
$ws.Cells.Item(73, 29).Value2 = 0
$ws.Cells.Item(73, 30).Value2 = 0

# Row 76
$ws.Cells.Item(76, 2).Value2 = 6839241
$ws.Cells.Item(76, 5).Value2 = 'Loughgall'
$ws.Cells.Item(76, 6).Value2 = 'Newry City'
$ws.Cells.Item(76, 7).Value2 = 3
$ws.Cells.Item(76, 8).Value2 = 1
$ws.Cells.Item(76, 9).Value2 = 0
$ws.Cells.Item(76, 10).Value2 = 0
$ws.Cells.Item(76, 11).Value2 = 'H'
$ws.Cells.Item(76, 12).Value2 = 1.7
$ws.Cells.Item(76, 13).Value2 = 3.6
$ws.Cells.Item(76, 14).Value2 = 4.5
$ws.Cells.Item(76, 15).Value2 = 1.833
$ws.Cells.Item(76, 16).Value2 = 3.6
$ws.Cells.Item(76, 17).Value2 = 3.8
$ws.Cells.Item(76, 18).Value2 = -0.5
$ws.Cells.Item(76, 19).Value2 = 1.825
$ws.Cells.Item(76, 20).Value2 = 1.975
$ws.Cells.Item(76, 21).Value2 = 2.5
$ws.Cells.Item(76, 22).Value2 = 1.8
$ws.Cells.Item(76, 23).Value2 = 2
$ws.Cells.Item(76, 24).Value2 = 0.833
$ws.Cells.Item(76, 25).Value2 = -1
$ws.Cells.Item(76, 26).Value2 = -1
$ws.Cells.Item(76, 27).Value2 = 0.825
$ws.Cells.Item(76, 28).Value2 = -1
$ws.Cells.Item(76, 29).Value2 = 0.8
$ws.Cells.Item(76, 30).Value2 = -1

# Row 77
$ws.Cells.Item(77, 2).Value2 = 6837584
$ws.Cells.Item(77, 5).Value2 = 'Glentoran'
$ws.Cells.Item(77, 6).Value2 = 'Glenavon'
$ws.Cells.Item(77, 7).Value2 = 3
$ws.Cells.Item(77, 8).Value2 = 1
$ws.Cells.Item(77, 9).Value2 = 3
$ws.Cells.Item(77, 10).Value2 = 0
$ws.Cells.Item(77, 11).Value2 = 'H'
$ws.Cells.Item(77, 12).Value2 = 1.3
$ws.Cells.Item(77, 13).Value2 = 5.25
$ws.Cells.Item(77, 14).Value2 = 7.5
$ws.Cells.Item(77, 15).Value2 = 1.3
$ws.Cells.Item(77, 16).Value2 = 5.5
$ws.Cells.Item(77, 17).Value2 = 7.5
$ws.Cells.Item(77, 18).Value2 = -1.5
$ws.Cells.Item(77, 19).Value2 = 1.85
$ws.Cells.Item(77, 20).Value2 = 1.95
$ws.Cells.Item(77, 21).Value2 = 3
$ws.Cells.Item(77, 22).Value2 = 1.85
$ws.Cells.Item(77, 23).Value2 = 1.95
$ws.Cells.Item(77, 24).Value2 = 0.3
$ws.Cells.Item(77, 25).Value2 = -1
$ws.Cells.Item(77, 26).Value2 = -1
$ws.Cells.Item(77, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(77, 28).Value2 = -1
$ws.Cells.Item(77, 29).Value2 = 0.8500000000000001
$ws.Cells.Item(77, 30).Value2 = -1

# Row 78
$ws.Cells.Item(78, 2).Value2 = 6840456
$ws.Cells.Item(78, 5).Value2 = 'Dungannon Swifts'
$ws.Cells.Item(78, 6).Value2 = 'Crusaders'
$ws.Cells.Item(78, 7).Value2 = 1
$ws.Cells.Item(78, 8).Value2 = 4
$ws.Cells.Item(78, 9).Value2 = 1
$ws.Cells.Item(78, 10).Value2 = 0
$ws.Cells.Item(78, 11).Value2 = 'A'
$ws.Cells.Item(78, 12).Value2 = 6
$ws.Cells.Item(78, 13).Value2 = 4.75
$ws.Cells.Item(78, 14).Value2 = 1.4
$ws.Cells.Item(78, 15).Value2 = 5.5
$ws.Cells.Item(78, 16).Value2 = 4.75
$ws.Cells.Item(78, 17).Value2 = 1.444
$ws.Cells.Item(78, 18).Value2 = 1.25
$ws.Cells.Item(78, 19).Value2 = 1.85
$ws.Cells.Item(78, 20).Value2 = 1.95
$ws.Cells.Item(78, 21).Value2 = 3
$ws.Cells.Item(78, 22).Value2 = 1.925
$ws.Cells.Item(78, 23).Value2 = 1.875
$ws.Cells.Item(78, 24).Value2 = -1
$ws.Cells.Item(78, 25).Value2 = -1
$ws.Cells.Item(78, 26).Value2 = 0.444
$ws.Cells.Item(78, 27).Value2 = -1
$ws.Cells.Item(78, 28).Value2 = 0.95
$ws.Cells.Item(78, 29).Value2 = 0.925
$ws.Cells.Item(78, 30).Value2 = -1

# Row 79
$ws.Cells.Item(79, 2).Value2 = 6840455
$ws.Cells.Item(79, 5).Value2 = 'Ballymena Utd'
$ws.Cells.Item(79, 6).Value2 = 'Coleraine'
$ws.Cells.Item(79, 7).Value2 = 3
$ws.Cells.Item(79, 8).Value2 = 1
$ws.Cells.Item(79, 9).Value2 = 0
$ws.Cells.Item(79, 10).Value2 = 0
$ws.Cells.Item(79, 11).Value2 = 'H'
$ws.Cells.Item(79, 12).Value2 = 4.75
$ws.Cells.Item(79, 13).Value2 = 3.8
$ws.Cells.Item(79, 14).Value2 = 1.615
$ws.Cells.Item(79, 15).Value2 = 4.5
$ws.Cells.Item(79, 16).Value2 = 3.8
$ws.Cells.Item(79, 17).Value2 = 1.615
$ws.Cells.Item(79, 18).Value2 = 0.75
$ws.Cells.Item(79, 19).Value2 = 2
$ws.Cells.Item(79, 20).Value2 = 1.8
$ws.Cells.Item(79, 21).Value2 = 2.75
$ws.Cells.Item(79, 22).Value2 = 2
$ws.Cells.Item(79, 23).Value2 = 1.8
$ws.Cells.Item(79, 24).Value2 = 3.5
$ws.Cells.Item(79, 25).Value2 = -1
$ws.Cells.Item(79, 26).Value2 = -1
$ws.Cells.Item(79, 27).Value2 = 1
$ws.Cells.Item(79, 28).Value2 = -1
$ws.Cells.Item(79, 29).Value2 = 1
$ws.Cells.Item(79, 30).Value2 = -1

# Row 151
$ws.Cells.Item(151, 2).Value2 = 6840490
$ws.Cells.Item(151, 5).Value2 = 'Dungannon Swifts'
$ws.Cells.Item(151, 6).Value2 = 'Larne FC'
$ws.Cells.Item(151, 7).Value2 = 0
$ws.Cells.Item(151, 8).Value2 = 2
$ws.Cells.Item(151, 9).Value2 = 0
$ws.Cells.Item(151, 10).Value2 = 2
$ws.Cells.Item(151, 11).Value2 = 'A'
$ws.Cells.Item(151, 12).Value2 = 5
$ws.Cells.Item(151, 13).Value2 = 4
$ws.Cells.Item(151, 14).Value2 = 1.5
$ws.Cells.Item(151, 15).Value2 = 6
$ws.Cells.Item(151, 16).Value2 = 4
$ws.Cells.Item(151, 17).Value2 = 1.444
$ws.Cells.Item(151, 18).Value2 = 1.25
$ws.Cells.Item(151, 19).Value2 = 1.8
$ws.Cells.Item(151, 20).Value2 = 2
$ws.Cells.Item(151, 21).Value2 = 2.75
$ws.Cells.Item(151, 22).Value2 = 2
$ws.Cells.Item(151, 23).Value2 = 1.8
$ws.Cells.Item(151, 24).Value2 = -1
$ws.Cells.Item(151, 25).Value2 = -1
$ws.Cells.Item(151, 26).Value2 = 0.444
$ws.Cells.Item(151, 27).Value2 = -1
$ws.Cells.Item(151, 28).Value2 = 1
$ws.Cells.Item(151, 29).Value2 = -1
$ws.Cells.Item(151, 30).Value2 = 0.8

# Row 152
$ws.Cells.Item(152, 2).Value2 = 6840294
$ws.Cells.Item(152, 5).Value2 = 'Linfield'
$ws.Cells.Item(152, 6).Value2 = 'Crusaders'
$ws.Cells.Item(152, 7).Value2 = 1
$ws.Cells.Item(152, 8).Value2 = 0
$ws.Cells.Item(152, 9).Value2 = 1
$ws.Cells.Item(152, 10).Value2 = 0
$ws.Cells.Item(152, 11).Value2 = 'H'
$ws.Cells.Item(152, 12).Value2 = 1.571
$ws.Cells.Item(152, 13).Value2 = 4
$ws.Cells.Item(152, 14).Value2 = 4.333
$ws.Cells.Item(152, 15).Value2 = 1.571
$ws.Cells.Item(152, 16).Value2 = 3.8
$ws.Cells.Item(152, 17).Value2 = 4.5
$ws.Cells.Item(152, 18).Value2 = -0.75
$ws.Cells.Item(152, 19).Value2 = 1.775
$ws.Cells.Item(152, 20).Value2 = 2.025
$ws.Cells.Item(152, 21).Value2 = 2.25
$ws.Cells.Item(152, 22).Value2 = 1.8
$ws.Cells.Item(152, 23).Value2 = 2
$ws.Cells.Item(152, 24).Value2 = 0.571
$ws.Cells.Item(152, 25).Value2 = -1
$ws.Cells.Item(152, 26).Value2 = -1
$ws.Cells.Item(152, 27).Value2 = 0.3875
$ws.Cells.Item(152, 28).Value2 = -0.5
$ws.Cells.Item(152, 29).Value2 = -1
$ws.Cells.Item(152, 30).Value2 = 1

# Row 153
$ws.Cells.Item(153, 2).Value2 = 6840489
$ws.Cells.Item(153, 5).Value2 = 'Coleraine'
$ws.Cells.Item(153, 6).Value2 = 'Carrick Rangers'
$ws.Cells.Item(153, 7).Value2 = 0
$ws.Cells.Item(153, 8).Value2 = 2
$ws.Cells.Item(153, 9).Value2 = 0
$ws.Cells.Item(153, 10).Value2 = 0
$ws.Cells.Item(153, 11).Value2 = 'A'
$ws.Cells.Item(153, 12).Value2 = 1.8
$ws.Cells.Item(153, 13).Value2 = 3.6
$ws.Cells.Item(153, 14).Value2 = 3.6
$ws.Cells.Item(153, 15).Value2 = 1.727
$ws.Cells.Item(153, 16).Value2 = 3.75
$ws.Cells.Item(153, 17).Value2 = 3.75
$ws.Cells.Item(153, 18).Value2 = -0.5
$ws.Cells.Item(153, 19).Value2 = 1.825
$ws.Cells.Item(153, 20).Value2 = 1.975
$ws.Cells.Item(153, 21).Value2 = 2.75
$ws.Cells.Item(153, 22).Value2 = 1.875
$ws.Cells.Item(153, 23).Value2 = 1.925
$ws.Cells.Item(153, 24).Value2 = -1
$ws.Cells.Item(153, 25).Value2 = -1
$ws.Cells.Item(153, 26).Value2 = 2.75
$ws.Cells.Item(153, 27).Value2 = -1
$ws.Cells.Item(153, 28).Value2 = 0.9750000000000001
$ws.Cells.Item(153, 29).Value2 = -1
$ws.Cells.Item(153, 30).Value2 = 0.925

# Row 154
$ws.Cells.Item(154, 2).Value2 = 6840293
$ws.Cells.Item(154, 5).Value2 = 'Cliftonville'
$ws.Cells.Item(154, 6).Value2 = 'Ballymena Utd'
$ws.Cells.Item(154, 7).Value2 = 3
$ws.Cells.Item(154, 8).Value2 = 0
$ws.Cells.Item(154, 9).Value2 = 1
$ws.Cells.Item(154, 10).Value2 = 0
$ws.Cells.Item(154, 11).Value2 = 'H'
$ws.Cells.Item(154, 12).Value2 = 1.125
$ws.Cells.Item(154, 13).Value2 = 8
$ws.Cells.Item(154, 14).Value2 = 17
$ws.Cells.Item(154, 15).Value2 = 1.166
$ws.Cells.Item(154, 16).Value2 = 7
$ws.Cells.Item(154, 17).Value2 = 13
$ws.Cells.Item(154, 18).Value2 = -2
$ws.Cells.Item(154, 19).Value2 = 1.875
$ws.Cells.Item(154, 20).Value2 = 1.925
$ws.Cells.Item(154, 21).Value2 = 3
$ws.Cells.Item(154, 22).Value2 = 1.85
$ws.Cells.Item(154, 23).Value2 = 1.95
$ws.Cells.Item(154, 24).Value2 = 0.1659999999999999
$ws.Cells.Item(154, 25).Value2 = -1
$ws.Cells.Item(154, 26).Value2 = -1
$ws.Cells.Item(154, 27).Value2 = 0.875
$ws.Cells.Item(154, 28).Value2 = -1
$ws.Cells.Item(154, 29).Value2 = 0
$ws.Cells.Item(154, 30).Value2 = 0

# Row 185
$ws.Cells.Item(185, 2).Value2 = 6840301
$ws.Cells.Item(185, 5).Value2 = 'Carrick Rangers'
$ws.Cells.Item(185, 6).Value2 = 'Newry City'
$ws.Cells.Item(185, 7).Value2 = 0
$ws.Cells.Item(185, 8).Value2 = 1
$ws.Cells.Item(185, 9).Value2 = 0
$ws.Cells.Item(185, 10).Value2 = 0
$ws.Cells.Item(185, 11).Value2 = 'A'
$ws.Cells.Item(185, 12).Value2 = 1.444
$ws.Cells.Item(185, 13).Value2 = 4
$ws.Cells.Item(185, 14).Value2 = 5.75
$ws.Cells.Item(185, 15).Value2 = 1.45
$ws.Cells.Item(185, 16).Value2 = 4
$ws.Cells.Item(185, 17).Value2 = 5.75
$ws.Cells.Item(185, 18).Value2 = -1
$ws.Cells.Item(185, 19).Value2 = 1.775
$ws.Cells.Item(185, 20).Value2 = 2.025
$ws.Cells.Item(185, 21).Value2 = 3
$ws.Cells.Item(185, 22).Value2 = 1.95
$ws.Cells.Item(185, 23).Value2 = 1.85
$ws.Cells.Item(185, 24).Value2 = -1
$ws.Cells.Item(185, 25).Value2 = -1
$ws.Cells.Item(185, 26).Value2 = 4.75
$ws.Cells.Item(185, 27).Value2 = -1
$ws.Cells.Item(185, 28).Value2 = 1.025
$ws.Cells.Item(185, 29).Value2 = -1
$ws.Cells.Item(185, 30).Value2 = 0.8500000000000001

# Row 186
$ws.Cells.Item(186, 2).Value2 = 6840302
$ws.Cells.Item(186, 5).Value2 = 'Glentoran'
$ws.Cells.Item(186, 6).Value2 = 'Loughgall'
$ws.Cells.Item(186, 7).Value2 = 0
$ws.Cells.Item(186, 8).Value2 = 3
$ws.Cells.Item(186, 9).Value2 = 0
$ws.Cells.Item(186, 10).Value2 = 2
$ws.Cells.Item(186, 11).Value2 = 'A'
$ws.Cells.Item(186, 12).Value2 = 1.4
$ws.Cells.Item(186, 13).Value2 = 4.333
$ws.Cells.Item(186, 14).Value2 = 6
$ws.Cells.Item(186, 15).Value2 = 1.222
$ws.Cells.Item(186, 16).Value2 = 5.75
$ws.Cells.Item(186, 17).Value2 = 8.5
$ws.Cells.Item(186, 18).Value2 = -1.75
$ws.Cells.Item(186, 19).Value2 = 1.85
$ws.Cells.Item(186, 20).Value2 = 1.95
$ws.Cells.Item(186, 21).Value2 = 3.5
$ws.Cells.Item(186, 22).Value2 = 2
$ws.Cells.Item(186, 23).Value2 = 1.8
$ws.Cells.Item(186, 24).Value2 = -1
$ws.Cells.Item(186, 25).Value2 = -1
$ws.Cells.Item(186, 26).Value2 = 7.5
$ws.Cells.Item(186, 27).Value2 = -1
$ws.Cells.Item(186, 28).Value2 = 0.95
$ws.Cells.Item(186, 29).Value2 = -1
$ws.Cells.Item(186, 30).Value2 = 0.8

# Row 188
$ws.Cells.Item(188, 2).Value2 = 6841446
$ws.Cells.Item(188, 5).Value2 = 'Larne FC'
$ws.Cells.Item(188, 6).Value2 = 'Glenavon'
$ws.Cells.Item(188, 7).Value2 = 6
$ws.Cells.Item(188, 8).Value2 = 1
$ws.Cells.Item(188, 9).Value2 = 5
$ws.Cells.Item(188, 10).Value2 = 1
$ws.Cells.Item(188, 11).Value2 = 'H'
$ws.Cells.Item(188, 12).Value2 = 1.333
$ws.Cells.Item(188, 13).Value2 = 5
$ws.Cells.Item(188, 14).Value2 = 8
$ws.Cells.Item(188, 15).Value2 = 1.4
$ws.Cells.Item(188, 16).Value2 = 4.5
$ws.Cells.Item(188, 17).Value2 = 7.5
$ws.Cells.Item(188, 18).Value2 = -1.25
$ws.Cells.Item(188, 19).Value2 = 1.825
$ws.Cells.Item(188, 20).Value2 = 1.975
$ws.Cells.Item(188, 21).Value2 = 2.75
$ws.Cells.Item(188, 22).Value2 = 1.825
$ws.Cells.Item(188, 23).Value2 = 1.975
$ws.Cells.Item(188, 24).Value2 = 0.3999999999999999
$ws.Cells.Item(188, 25).Value2 = -1
$ws.Cells.Item(188, 26).Value2 = -1
$ws.Cells.Item(188, 27).Value2 = 0.825
$ws.Cells.Item(188, 28).Value2 = -1
$ws.Cells.Item(188, 29).Value2 = 0.825
$ws.Cells.Item(188, 30).Value2 = -1

# Row 189
$ws.Cells.Item(189, 2).Value2 = 6840961
$ws.Cells.Item(189, 5).Value2 = 'Loughgall'
$ws.Cells.Item(189, 6).Value2 = 'Cliftonville'
$ws.Cells.Item(189, 7).Value2 = 2
$ws.Cells.Item(189, 8).Value2 = 3
$ws.Cells.Item(189, 9).Value2 = 1
$ws.Cells.Item(189, 10).Value2 = 1
$ws.Cells.Item(189, 11).Value2 = 'A'
$ws.Cells.Item(189, 12).Value2 = 7
$ws.Cells.Item(189, 13).Value2 = 4.5
$ws.Cells.Item(189, 14).Value2 = 1.4
$ws.Cells.Item(189, 15).Value2 = 4.2
$ws.Cells.Item(189, 16).Value2 = 4
$ws.Cells.Item(189, 17).Value2 = 1.65
$ws.Cells.Item(189, 18).Value2 = 0.75
$ws.Cells.Item(189, 19).Value2 = 1.95
$ws.Cells.Item(189, 20).Value2 = 1.85
$ws.Cells.Item(189, 21).Value2 = 3
$ws.Cells.Item(189, 22).Value2 = 1.925
$ws.Cells.Item(189, 23).Value2 = 1.875
$ws.Cells.Item(189, 24).Value2 = -1
$ws.Cells.Item(189, 25).Value2 = -1
$ws.Cells.Item(189, 26).Value2 = 0.6499999999999999
$ws.Cells.Item(189, 27).Value2 = -0.5
$ws.Cells.Item(189, 28).Value2 = 0.425
$ws.Cells.Item(189, 29).Value2 = 0.925
$ws.Cells.Item(189, 30).Value2 = -1

# Row 190
$ws.Cells.Item(190, 2).Value2 = 6840960
$ws.Cells.Item(190, 5).Value2 = 'Linfield'
$ws.Cells.Item(190, 6).Value2 = 'Newry City'
$ws.Cells.Item(190, 7).Value2 = 6
$ws.Cells.Item(190, 8).Value2 = 0
$ws.Cells.Item(190, 9).Value2 = 4
$ws.Cells.Item(190, 10).Value2 = 0
$ws.Cells.Item(190, 11).Value2 = 'H'
$ws.Cells.Item(190, 12).Value2 = 1.142
$ws.Cells.Item(190, 13).Value2 = 8.5
$ws.Cells.Item(190, 14).Value2 = 13
$ws.Cells.Item(190, 15).Value2 = 1.125
$ws.Cells.Item(190, 16).Value2 = 8.5
$ws.Cells.Item(190, 17).Value2 = 13
$ws.Cells.Item(190, 18).Value2 = -2.25
$ws.Cells.Item(190, 19).Value2 = 1.8
$ws.Cells.Item(190, 20).Value2 = 2
$ws.Cells.Item(190, 21).Value2 = 3.75
$ws.Cells.Item(190, 22).Value2 = 1.95
$ws.Cells.Item(190, 23).Value2 = 1.85
$ws.Cells.Item(190, 24).Value2 = 0.125
$ws.Cells.Item(190, 25).Value2 = -1
$ws.Cells.Item(190, 26).Value2 = -1
$ws.Cells.Item(190, 27).Value2 = 0.8
$ws.Cells.Item(190, 28).Value2 = -1
$ws.Cells.Item(190, 29).Value2 = 0.95
$ws.Cells.Item(190, 30).Value2 = -1

# Row 191
$ws.Cells.Item(191, 2).Value2 = 6840959
$ws.Cells.Item(191, 5).Value2 = 'Dungannon Swifts'
$ws.Cells.Item(191, 6).Value2 = 'Glentoran'
$ws.Cells.Item(191, 7).Value2 = 0
$ws.Cells.Item(191, 8).Value2 = 3
$ws.Cells.Item(191, 9).Value2 = 0
$ws.Cells.Item(191, 10).Value2 = 2
$ws.Cells.Item(191, 11).Value2 = 'A'
$ws.Cells.Item(191, 12).Value2 = 3.6
$ws.Cells.Item(191, 13).Value2 = 4
$ws.Cells.Item(191, 14).Value2 = 1.8
$ws.Cells.Item(191, 15).Value2 = 3.3
$ws.Cells.Item(191, 16).Value2 = 3.8
$ws.Cells.Item(191, 17).Value2 = 1.909
$ws.Cells.Item(191, 18).Value2 = 0.5
$ws.Cells.Item(191, 19).Value2 = 1.85
$ws.Cells.Item(191, 20).Value2 = 1.95
$ws.Cells.Item(191, 21).Value2 = 3
$ws.Cells.Item(191, 22).Value2 = 1.9
$ws.Cells.Item(191, 23).Value2 = 1.9
$ws.Cells.Item(191, 24).Value2 = -1
$ws.Cells.Item(191, 25).Value2 = -1
$ws.Cells.Item(191, 26).Value2 = 0.909
$ws.Cells.Item(191, 27).Value2 = -1
$ws.Cells.Item(191, 28).Value2 = 0.95
$ws.Cells.Item(191, 29).Value2 = 0
$ws.Cells.Item(191, 30).Value2 = 0

# Row 193
$ws.Cells.Item(193, 2).Value2 = 6840958
$ws.Cells.Item(193, 5).Value2 = 'Ballymena Utd'
$ws.Cells.Item(193, 6).Value2 = 'Carrick Rangers'
$ws.Cells.Item(193, 7).Value2 = 0
$ws.Cells.Item(193, 8).Value2 = 2
$ws.Cells.Item(193, 9).Value2 = 0
$ws.Cells.Item(193, 10).Value2 = 1
$ws.Cells.Item(193, 11).Value2 = 'A'
$ws.Cells.Item(193, 12).Value2 = 3
$ws.Cells.Item(193, 13).Value2 = 3.4
$ws.Cells.Item(193, 14).Value2 = 2.2
$ws.Cells.Item(193, 15).Value2 = 3.6
$ws.Cells.Item(193, 16).Value2 = 3.5
$ws.Cells.Item(193, 17).Value2 = 1.95
$ws.Cells.Item(193, 18).Value2 = 0.5
$ws.Cells.Item(193, 19).Value2 = 1.825
$ws.Cells.Item(193, 20).Value2 = 1.975
$ws.Cells.Item(193, 21).Value2 = 2.5
$ws.Cells.Item(193, 22).Value2 = 1.975
$ws.Cells.Item(193, 23).Value2 = 1.825
$ws.Cells.Item(193, 24).Value2 = -1
$ws.Cells.Item(193, 25).Value2 = -1
$ws.Cells.Item(193, 26).Value2 = 0.95
$ws.Cells.Item(193, 27).Value2 = -1
$ws.Cells.Item(193, 28).Value2 = 0.9750000000000001
$ws.Cells.Item(193, 29).Value2 = -1
$ws.Cells.Item(193, 30).Value2 = 0.825

# Row 196
$ws.Cells.Item(196, 2).Value2 = 6841450
$ws.Cells.Item(196, 5).Value2 = 'Glenavon'
$ws.Cells.Item(196, 6).Value2 = 'Loughgall'
$ws.Cells.Item(196, 7).Value2 = 1
$ws.Cells.Item(196, 8).Value2 = 2
$ws.Cells.Item(196, 9).Value2 = 0
$ws.Cells.Item(196, 10).Value2 = 2
$ws.Cells.Item(196, 11).Value2 = 'A'
$ws.Cells.Item(196, 12).Value2 = 1.8
$ws.Cells.Item(196, 13).Value2 = 3.75
$ws.Cells.Item(196, 14).Value2 = 3.4
$ws.Cells.Item(196, 15).Value2 = 1.8
$ws.Cells.Item(196, 16).Value2 = 3.75
$ws.Cells.Item(196, 17).Value2 = 3.5
$ws.Cells.Item(196, 18).Value2 = -0.5
$ws.Cells.Item(196, 19).Value2 = 1.8
$ws.Cells.Item(196, 20).Value2 = 2
$ws.Cells.Item(196, 21).Value2 = 2.75
$ws.Cells.Item(196, 22).Value2 = 1.8
$ws.Cells.Item(196, 23).Value2 = 2
$ws.Cells.Item(196, 24).Value2 = -1
$ws.Cells.Item(196, 25).Value2 = -1
$ws.Cells.Item(196, 26).Value2 = 2.5
$ws.Cells.Item(196, 27).Value2 = -1
$ws.Cells.Item(196, 28).Value2 = 1
$ws.Cells.Item(196, 29).Value2 = 0.4
$ws.Cells.Item(196, 30).Value2 = -0.5

# Row 197
$ws.Cells.Item(197, 2).Value2 = 6841449
$ws.Cells.Item(197, 5).Value2 = 'Coleraine'
$ws.Cells.Item(197, 6).Value2 = 'Dungannon Swifts'
$ws.Cells.Item(197, 7).Value2 = 1
$ws.Cells.Item(197, 8).Value2 = 1
$ws.Cells.Item(197, 9).Value2 = 1
$ws.Cells.Item(197, 10).Value2 = 0
$ws.Cells.Item(197, 11).Value2 = 'D'
$ws.Cells.Item(197, 12).Value2 = 1.85
$ws.Cells.Item(197, 13).Value2 = 4
$ws.Cells.Item(197, 14).Value2 = 3.1
$ws.Cells.Item(197, 15).Value2 = 1.8
$ws.Cells.Item(197, 16).Value2 = 4
$ws.Cells.Item(197, 17).Value2 = 3.2
$ws.Cells.Item(197, 18).Value2 = -0.5
$ws.Cells.Item(197, 19).Value2 = 1.875
$ws.Cells.Item(197, 20).Value2 = 1.925
$ws.Cells.Item(197, 21).Value2 = 3
$ws.Cells.Item(197, 22).Value2 = 1.95
$ws.Cells.Item(197, 23).Value2 = 1.85
$ws.Cells.Item(197, 24).Value2 = -1
$ws.Cells.Item(197, 25).Value2 = 3
$ws.Cells.Item(197, 26).Value2 = -1
$ws.Cells.Item(197, 27).Value2 = -1
$ws.Cells.Item(197, 28).Value2 = 0.925
$ws.Cells.Item(197, 29).Value2 = -1
$ws.Cells.Item(197, 30).Value2 = 0.8500000000000001

# Row 215
$ws.Cells.Item(215, 2).Value2 = 8048805
$ws.Cells.Item(215, 5).Value2 = 'Dungannon Swifts'
$ws.Cells.Item(215, 6).Value2 = 'Loughgall'
$ws.Cells.Item(215, 7).Value2 = 2
$ws.Cells.Item(215, 8).Value2 = 0
$ws.Cells.Item(215, 9).Value2 = 0
$ws.Cells.Item(215, 10).Value2 = 0
$ws.Cells.Item(215, 11).Value2 = 'H'
$ws.Cells.Item(215, 12).Value2 = 1.615
$ws.Cells.Item(215, 13).Value2 = 4
$ws.Cells.Item(215, 14).Value2 = 4.333
$ws.Cells.Item(215, 15).Value2 = 1.85
$ws.Cells.Item(215, 16).Value2 = 3.6
$ws.Cells.Item(215, 17).Value2 = 3.5
$ws.Cells.Item(215, 18).Value2 = -0.5
$ws.Cells.Item(215, 19).Value2 = 1.875
$ws.Cells.Item(215, 20).Value2 = 1.925
$ws.Cells.Item(215, 21).Value2 = 3
$ws.Cells.Item(215, 22).Value2 = 1.9
$ws.Cells.Item(215, 23).Value2 = 1.9
$ws.Cells.Item(215, 24).Value2 = 0.8500000000000001
$ws.Cells.Item(215, 25).Value2 = -1
$ws.Cells.Item(215, 26).Value2 = -1
$ws.Cells.Item(215, 27).Value2 = 0.875
$ws.Cells.Item(215, 28).Value2 = -1
$ws.Cells.Item(215, 29).Value2 = -1
$ws.Cells.Item(215, 30).Value2 = 0.8999999999999999

# Row 216
$ws.Cells.Item(216, 2).Value2 = 8048806
$ws.Cells.Item(216, 5).Value2 = 'Newry City'
$ws.Cells.Item(216, 6).Value2 = 'Carrick Rangers'
$ws.Cells.Item(216, 7).Value2 = 1
$ws.Cells.Item(216, 8).Value2 = 3
$ws.Cells.Item(216, 9).Value2 = 1
$ws.Cells.Item(216, 10).Value2 = 3
$ws.Cells.Item(216, 11).Value2 = 'A'
$ws.Cells.Item(216, 12).Value2 = 4
$ws.Cells.Item(216, 13).Value2 = 3.8
$ws.Cells.Item(216, 14).Value2 = 1.727
$ws.Cells.Item(216, 15).Value2 = 4.2
$ws.Cells.Item(216, 16).Value2 = 3.8
$ws.Cells.Item(216, 17).Value2 = 1.7
$ws.Cells.Item(216, 18).Value2 = 0.75
$ws.Cells.Item(216, 19).Value2 = 1.9
$ws.Cells.Item(216, 20).Value2 = 1.9
$ws.Cells.Item(216, 21).Value2 = 2.75
$ws.Cells.Item(216, 22).Value2 = 1.825
$ws.Cells.Item(216, 23).Value2 = 1.975
$ws.Cells.Item(216, 24).Value2 = -1
$ws.Cells.Item(216, 25).Value2 = -1
$ws.Cells.Item(216, 26).Value2 = 0.7
$ws.Cells.Item(216, 27).Value2 = -1
$ws.Cells.Item(216, 28).Value2 = 0.8999999999999999
$ws.Cells.Item(216, 29).Value2 = 0.825
$ws.Cells.Item(216, 30).Value2 = -1

# Row 217
$ws.Cells.Item(217, 2).Value2 = 8048813
$ws.Cells.Item(217, 5).Value2 = 'Glentoran'
$ws.Cells.Item(217, 6).Value2 = 'Larne FC'
$ws.Cells.Item(217, 7).Value2 = 1
$ws.Cells.Item(217, 8).Value2 = 2
$ws.Cells.Item(217, 9).Value2 = 0
$ws.Cells.Item(217, 10).Value2 = 1
$ws.Cells.Item(217, 11).Value2 = 'A'
$ws.Cells.Item(217, 12).Value2 = 3
$ws.Cells.Item(217, 13).Value2 = 3.6
$ws.Cells.Item(217, 14).Value2 = 2.05
$ws.Cells.Item(217, 15).Value2 = 4.5
$ws.Cells.Item(217, 16).Value2 = 4
$ws.Cells.Item(217, 17).Value2 = 1.615
$ws.Cells.Item(217, 18).Value2 = 0.75
$ws.Cells.Item(217, 19).Value2 = 2.025
$ws.Cells.Item(217, 20).Value2 = 1.775
$ws.Cells.Item(217, 21).Value2 = 2.75
$ws.Cells.Item(217, 22).Value2 = 1.85
$ws.Cells.Item(217, 23).Value2 = 1.95
$ws.Cells.Item(217, 24).Value2 = -1
$ws.Cells.Item(217, 25).Value2 = -1
$ws.Cells.Item(217, 26).Value2 = 0.615
$ws.Cells.Item(217, 27).Value2 = -0.5
$ws.Cells.Item(217, 28).Value2 = 0.3875
$ws.Cells.Item(217, 29).Value2 = 0.425
$ws.Cells.Item(217, 30).Value2 = -0.5
